# Edit script: update cluster outbreak data
# New data set (89 rows) to replace existing data (84 rows) starting at row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("3535 Opal Meadow Heights Aged Care Community Meadow Heights", 22),
    @("574 Plummer Street Building B1 & B2 Port Melbourne", 5),
    @("7 Chefs Fawkner", 6),
    @("A1 Bakery Brunswick", 6),
    @("Acquire BPO Southbank", 15),
    @("Al Haj Halal Meats Glenroy", 52),
    @("Al-Taqwa College Truganina", 14),
    @("Amiga Montessori Craigieburn", 25),
    @("Best&Less Fountain Gate Narre Warren", 6),
    @("Broadmeadows Medical Centre Broadmeadows", 5),
    @("Budget Car and Truck Rentals Campbellfield", 7),
    @("Campbellfield Heights Primary School Campbellfield", 6),
    @("Can Panel Cambellfield", 6),
    @("Cannie Road Construction Site Cannie", 7),
    @("Caroline Springs Police Station", 7),
    @("Cedars Medical Clinic Coburg", 49),
    @("Central Plaster 236 & 259 Jasper Road McKinnon", 6),
    @("City of Hobsons Bay Community", 10),
    @("City of Moreland Community", 7),
    @("City of Wyndham Community", 5),
    @("Classy Cabinets and Kitchens Craigieburn", 13),
    @("Coles Barkly Square Brunswick August", 5),
    @("Coles Broadmeadows Central Shopping Centre", 10),
    @("Coles Campbellfield Plaza Campbellfield", 9),
    @("Coles Coburg North Village", 14),
    @("Coles Coburg North Village Aug", 7),
    @("Coles Greenvale Shopping Centre", 6),
    @("Coles Pakenham Place Shopping Centre", 7),
    @("Coles Roxburgh Village Roxburgh Park", 15),
    @("Community Kids Meadow Heights", 21),
    @("Construction Site Olea Apartment Caulfield North", 5),
    @("Costco Wholesale Epping", 22),
    @("Crossroads Logistics Sunshine North", 5),
    @("Croydon Orthodontics", 6),
    @("Elite Smart Community Care Campbellfield", 5),
    @("Ernst and Young 8 Exhibition Street Melbourne", 5),
    @("Fitzroy Community School Fitzroy North", 57),
    @("Glenroy West Primary School", 7),
    @("Gloria Jeans Coffees Broadmeadows Central", 5),
    @("Health Care Providers Association South Melbourne", 16),
    @("IGA Meadow Heights Shopping Centre Meadow Heights", 5),
    @("Ilim College Glenroy Campus Hadfield", 22),
    @("Ilim College Kiewa Campus Boys Secondary Dallas", 5),
    @("Ilim Learning Sanctuary Glenroy", 18),
    @("Industrial Galvanizers Valmont Coatings Campbellfield", 15),
    @("Islamic College of Melbourne Tarneit", 9),
    @("KFC Fawkner", 14),
    @("Kasr Sweets Coolaroo", 11),
    @("Kids House Early Learning Cheltenham", 13),
    @("Learning Nest Early Learning Centre Meadow Heights", 10),
    @("Level Crossing Removal Project Lilydale Construction Site John Street", 9),
    @("Linfox Somerton National Distribution Centre Somerton", 7),
    @("Malvern Health and Fitness Clinic Malvern", 6),
    @("Melbourne Metropolitan Remand Centre Ravenhall", 9),
    @("Melbourne Truck Repairs Campbellfield", 6),
    @("Melbourne West Police Station Docklands", 6),
    @("Melbourne Youth Justice Centre Parkville", 5),
    @("Mercy Hospital for Women Heidelberg", 6),
    @("Mernda YMCA Early Learning Centre Mernda", 5),
    @("Montessori Beginnings Greenvale", 6),
    @("MyCentre Childcare Broadmeadows", 15),
    @("Newbury Child and Community Centre Craigieburn", 10),
    @("Newport Football Club Altona North", 5),
    @("Newport Gardens Early Years Centre Newport", 5),
    @("Nido Early School Moonee Ponds", 18),
    @("Nino Early Learning Adventures Lalor", 5),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 48),
    @("OnQ Plumbing and Excavations Craigieburn", 11),
    @("Oporto Coolaroo", 15),
    @("Paisley Park Early Learning Centre Bundoora", 10),
    @("Panorama Construction Site Whitehorse Rd Box Hill", 40),
    @("People First Healthcare Home Residence Disability Support Taylors Lakes", 5),
    @("Salta Drive Construction Site Rangedale Drainage Altona North", 6),
    @("Serco Mill Park", 6),
    @("Southern Cross Station Crew Room Tissue Box Docklands", 5),
    @("St Vincents Hospital Emergency Department Melbourne", 5),
    @("Tek Foods Somerton", 10),
    @("The Homestead Child and Family Centre Roxburgh Park", 18),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 6),
    @("Tip Top Warehouse Dandenong", 10),
    @("Total Window Concepts Hoppers Crossing", 6),
    @("Tunstall Fresh Tunstall Square Shopping Centre Doncaster East", 5),
    @("Unilodge College Square Student Accommodation 570 Lygon Street Carlton", 10),
    @("Victorian Civil and Administrative Tribunal Melbourne", 5),
    @("Western Health Sunshine Hospital Emergency Department", 5),
    @("Who is Bunker Spreckels Cafe Elwood", 5),
    @("Woodlands Long Day Care and Kindergarten Roxburgh Park", 5),
    @("Woolworths Greenvale Lakes Roxburgh Park", 13),
    @("Yara Childcare Centre Truganina", 26)
)

$numRows = $data.Count

# Clear any existing data rows below the header (old sheet had rows 2..85)
$ws.Range("A2:B1000").Clear()

# Write header (unchanged, but ensure it's correct)
$ws.Cells.Item(1, 1).Value = "Cluster name"
$ws.Cells.Item(1, 2).Value = "Active cases"

# Write the new rows starting at row 2
for ($i = 0; $i -lt $numRows; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
